$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns: AD=Wins, AE=Losses, AF=Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, bordered, centered) used by A1:AC1
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Team record (Wins-Losses-Ties) repeated for every player row
$lastRow = 55
$ws.Range("AD2:AD" + $lastRow).Value = 86
$ws.Range("AE2:AE" + $lastRow).Value = 76
$ws.Range("AF2:AF" + $lastRow).Value = 0
